$d = $word.ActiveDocument

# Desired merged text per paragraph style (runs in these paragraphs are
# currently split word-by-word / space-by-space; collapse each into a
# single run holding the full paragraph text).
$targets = @{
    "Title"    = "Questions: The scalar product"
    "Author"   = "Ritwik Anand"
    "Abstract" = "A selection of questions for the study guide on the scalar product"
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    if ($targets.ContainsKey($styleName)) {
        $text = $targets[$styleName]
        $start = $p.Range.Start
        $end = $p.Range.End - 1

        # Use a temporary sentinel change first so the runtime always
        # treats the second assignment as a genuine text change, which
        # forces it to collapse the paragraph's many runs into one.
        $rTemp = $d.Range($start, $end)
        $rTemp.Text = $text + [char]1

        $p2 = $d.Paragraphs.Item($i)
        $rFinal = $d.Range($p2.Range.Start, $p2.Range.End - 1)
        $rFinal.Text = $text
    }
}
